$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "Scanner" to "Session"
$ws.Name = "Session"

# Delete row 2 (the extra data row) so only the header row remains
$ws.Rows.Item(2).Delete()
